# Apply the "Create Client & Corrected Excels" edits to the workbook.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Summary sheet: correct a few figures in row 4
# ------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 150
$wsSummary.Range("B4").Value = 50
$wsSummary.Range("E4").Value = 100

# ------------------------------------------------------------------
# Repayment schedule sheet: correct fee/outstanding figures
# ------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("I3").Value = 50
$wsRepay.Range("K3").Value = 937.72
$wsRepay.Range("L3").Value = 937.72

$wsRepay.Range("I5").Value = 0
$wsRepay.Range("K5").Value = 887.72
$wsRepay.Range("L5").Value = 887.72

$wsRepay.Range("I6").Value = 100
$wsRepay.Range("K6").Value = 987.72
$wsRepay.Range("P6").Value = 987.72

# ------------------------------------------------------------------
# Transactions sheet: correct transaction ids / fee figures
# ------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 47
$wsTrans.Range("E2").Value = 887.72
$wsTrans.Range("H2").Value = 0

$wsTrans.Range("A3").Value = 44

$wsTrans.Range("A4").Value = 42
$wsTrans.Range("E4").Value = 937.72
$wsTrans.Range("H4").Value = 50

$wsTrans.Range("A5").Value = 38

# ------------------------------------------------------------------
# Update selections / active cells on each sheet to match the
# author's final cursor positions.
# ------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("E19").Select()

$wsSummary.Range("D26").Select()

$wsRepay.Range("F22").Select()

$wsTrans.Range("E2:G2").Select()

# ------------------------------------------------------------------
# Switch the active/selected sheet from Transactions to Acc_Repayment1
# (this also updates workbookView's activeTab and each sheet's
# tabSelected flag).
# ------------------------------------------------------------------
$wsAccRepayment1 = $wb.Worksheets.Item("Acc_Repayment1")
$wsAccRepayment1.Activate()
$wsAccRepayment1.Range("G3").Select()
